# "origin environment test setup"
# The log sheet's A2 cell tracks the most recent build. Two new builds ran
# since the last recorded one (3feea22dad @ 2020-07-24 15:52); append them
# in order so the log ends up pointing at the latest build,
# 0c48d8e2bd @ 2020-08-13 12:02.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "236bbbf490 built at 2020-08-12 17:11`n"
$ws.Range("A2").Value = "0c48d8e2bd built at 2020-08-13 12:02`n"
